$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (shifts existing rows 9+ down by one)
$ws.Rows.Item(9).Insert()

# Populate new row 9: un_franzosa_ControlvsCD_ConvCD
$ws.Range("A9").Value = "un_franzosa_ControlvsCD_ConvCD"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.4
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6
$ws.Range("H9").Value = 0.6

# Insert a new row at row 14 (shifts existing rows 14+ down by one)
$ws.Rows.Item(14).Insert()

# Populate new row 14: un_franzosa_ControlvsUC_ConvUC
$ws.Range("A14").Value = "un_franzosa_ControlvsUC_ConvUC"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0.4
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6
$ws.Range("H14").Value = 0.6

# Populate new appended row 25: nf_wang_egfr
$ws.Range("A25").Value = "nf_wang_egfr"
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0.6
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 0.4
$ws.Range("H25").Value = 0.4

# Populate new appended row 26: nf_wang_studygroup
$ws.Range("A26").Value = "nf_wang_studygroup"
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0.6
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 0.4
$ws.Range("H26").Value = 0.4

# Populate new appended row 27: nf_wang_urea
$ws.Range("A27").Value = "nf_wang_urea"
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0.8
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 0.2
$ws.Range("H27").Value = 0.2
